# Apply "modificacion de usuario" edit:
#  - Sheet "Registros": correct Martes (row 3) punch times + hour totals,
#    and append four new daily rows (Miercoles..Sabado) for legajo 52.
#  - Sheet "Totalizado": refresh the aggregated totals row to match the
#    updated detail rows.

$wb = $excel.ActiveWorkbook
$registros = $wb.Worksheets.Item("Registros")
$totalizado = $wb.Worksheets.Item("Totalizado")

$nombre = "SIRAGUSA JAVIER CARLOS S. "

function Set-Punch {
    param($ws, [int]$row, [string]$col, [double]$value)
    $cell = $ws.Range("$col$row")
    $cell.Value = $value
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

function Set-FechaDate {
    param($ws, [int]$row, [string]$col, [double]$value)
    $cell = $ws.Range("$col$row")
    $cell.Value = $value
    $cell.NumberFormat = "YYYY-MM-DD"
}

function Set-DayRow {
    param($ws, [int]$row, [int]$legajo, [string]$nombre, [string]$dia, [double]$fecha,
          [double]$e0, [double]$s0, [double]$e1, [double]$s1, [double]$e2, [double]$s2,
          [double]$e3, [double]$s3, [double]$e4, [double]$s4,
          [double]$hnorm, [double]$h50, [double]$h100)

    $ws.Range("A$row").Value = $legajo
    $ws.Range("B$row").Value = $nombre
    $ws.Range("C$row").Value = $dia
    Set-FechaDate $ws $row "D" $fecha
    Set-Punch $ws $row "E" $e0
    Set-Punch $ws $row "F" $s0
    Set-Punch $ws $row "G" $e1
    Set-Punch $ws $row "H" $s1
    Set-Punch $ws $row "I" $e2
    Set-Punch $ws $row "J" $s2
    Set-Punch $ws $row "K" $e3
    Set-Punch $ws $row "L" $s3
    Set-Punch $ws $row "M" $e4
    Set-Punch $ws $row "N" $s4
    $ws.Range("O$row").Value = $hnorm
    $ws.Range("P$row").Value = $h50
    $ws.Range("Q$row").Value = $h100
}

# --- Row 3 (Martes) gets corrected punch times + recomputed hour buckets ---
Set-DayRow $registros 3 52 $nombre "Martes" 44138 `
    44138.35 44138.52777777778 44138.55069444444 44138.74236111111 `
    44138 44138 44138 44138 44138 44138 `
    7.6 1 0

# --- New row 4 (Miercoles) ---
Set-DayRow $registros 4 52 $nombre "Miércoles" 44139 `
    44139.32777777778 44139.74236111111 44139 44139 `
    44139 44139 44139 44139 44139 44139 `
    8 1.12 0

# --- New row 5 (Jueves) ---
Set-DayRow $registros 5 52 $nombre "Jueves" 44140 `
    44140.32777777778 44140.71527777778 44140 44140 `
    44140 44140 44140 44140 44140 44140 `
    8 0.47 0

# --- New row 6 (Viernes) ---
Set-DayRow $registros 6 52 $nombre "Viernes" 44141 `
    44141.32569444443 44141.52708333332 44141.55 44141.74236111111 `
    44141 44141 44141 44141 44141 44141 `
    8 1.17 0

# --- New row 7 (Sabado) ---
Set-DayRow $registros 7 52 $nombre "Sábado" 44142 `
    44142.32569444443 44142.50069444445 44142 44142 `
    44142 44142 44142 44142 44142 44142 `
    0 0 0

# --- Refresh the "Totalizado" aggregate row to match the new detail rows ---
$totalizado.Range("C2").Value = 39.6
$totalizado.Range("D2").Value = 4.88
$totalizado.Range("E2").Value = 0
$totalizado.Range("F2").Value = 6
